# Auto-generated Excel COM-interop script to apply Asura_Profits.xlsx edits
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H15").Value = 1259.8833
$ws.Range("I15").Value = 1259.8833
$ws.Range("K15").Value = 3779.6499
$ws.Range("M15").Value = -3610.6499
$ws.Range("H74").Value = 4549.9375
$ws.Range("I74").Value = 4286.143
$ws.Range("J74").Value = 4755.1113
$ws.Range("K74").Value = 4286.143
$ws.Range("L74").Value = 4755.1113
$ws.Range("M74").Value = -3350.143
$ws.Range("N74").Value = -6627.1113
$ws.Range("H77").Value = 4549.9375
$ws.Range("I77").Value = 4286.143
$ws.Range("J77").Value = 4755.1113
$ws.Range("K77").Value = 21430.715
$ws.Range("L77").Value = 23775.5565
$ws.Range("M77").Value = -16750.715
$ws.Range("N77").Value = -33135.5565
$ws.Range("H88").Value = 2733.3333
$ws.Range("J88").Value = 3028.5715
$ws.Range("L88").Value = 3028.5715
$ws.Range("N88").Value = -3840.5715
$ws.Range("H91").Value = 2733.3333
$ws.Range("J91").Value = 3028.5715
$ws.Range("L91").Value = 3028.5715
$ws.Range("N91").Value = -5836.5715
$ws.Range("H113").Value = 2924.7058
$ws.Range("I113").Value = 2347.1428
$ws.Range("J113").Value = 3329
$ws.Range("K113").Value = 2347.1428
$ws.Range("L113").Value = 3329
$ws.Range("M113").Value = 906.8571999999999
$ws.Range("N113").Value = -9837
$ws.Range("H129").Value = 1160.1471
$ws.Range("J129").Value = 1644
$ws.Range("L129").Value = 4932
$ws.Range("N129").Value = -14932
$ws.Range("H131").Value = 6508.1934
$ws.Range("I131").Value = 1004.3333
$ws.Range("J131").Value = 8759.772000000001
$ws.Range("K131").Value = 3012.9999
$ws.Range("L131").Value = 26279.316
$ws.Range("M131").Value = 2027.0001
$ws.Range("N131").Value = -36359.31600000001
$ws.Range("H137").Value = 1438.32
$ws.Range("I137").Value = 1318.6842
$ws.Range("J137").Value = 1817.1666
$ws.Range("K137").Value = 3956.0526
$ws.Range("L137").Value = 5451.4998
$ws.Range("M137").Value = -1406.0526
$ws.Range("N137").Value = -10551.4998
$ws.Range("H138").Value = 2414.7917
$ws.Range("I138").Value = 1319.9048
$ws.Range("J138").Value = 3947.6333
$ws.Range("K138").Value = 3959.7144
$ws.Range("L138").Value = 11842.8999
$ws.Range("M138").Value = 1180.2856
$ws.Range("N138").Value = -22122.8999
$ws.Range("H141").Value = 7101.364
$ws.Range("I141").Value = 2475.2632
$ws.Range("J141").Value = 36400
$ws.Range("K141").Value = 7425.7896
$ws.Range("L141").Value = 109200
$ws.Range("M141").Value = -2245.7896
$ws.Range("N141").Value = -119560
$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H2").Value = 24877.143
$ws.Range("I2").Value = 867.36365
$ws.Range("J2").Value = 51287.9
$ws.Range("K2").Value = 867.36365
$ws.Range("L2").Value = 51287.9
$ws.Range("M2").Value = -754.36365
$ws.Range("N2").Value = -51513.9
$ws.Range("H61").Value = 2561.125
$ws.Range("I61").Value = 2095.8572
$ws.Range("J61").Value = 3449.3635
$ws.Range("K61").Value = 2095.8572
$ws.Range("L61").Value = 3449.3635
$ws.Range("M61").Value = -1883.8572
$ws.Range("N61").Value = -3873.3635
$ws.Range("H74").Value = 1040.8485
$ws.Range("I74").Value = 962.8461
$ws.Range("J74").Value = 1330.5714
$ws.Range("K74").Value = 962.8461
$ws.Range("L74").Value = 1330.5714
$ws.Range("M74").Value = -88.84609999999998
$ws.Range("N74").Value = -3078.5714
$ws.Range("H77").Value = 1040.8485
$ws.Range("I77").Value = 962.8461
$ws.Range("J77").Value = 1330.5714
$ws.Range("K77").Value = 4814.2305
$ws.Range("L77").Value = 6652.857
$ws.Range("M77").Value = -446.2304999999997
$ws.Range("N77").Value = -15388.857
$ws.Range("H116").Value = 24877.143
$ws.Range("I116").Value = 867.36365
$ws.Range("J116").Value = 51287.9
$ws.Range("K116").Value = 867.36365
$ws.Range("L116").Value = 51287.9
$ws.Range("M116").Value = 1426.63635
$ws.Range("N116").Value = -55875.9
$ws.Range("H122").Value = 5644.2104
$ws.Range("I122").Value = 7462.375
$ws.Range("J122").Value = 4321.909
$ws.Range("K122").Value = 22387.125
$ws.Range("L122").Value = 12965.727
$ws.Range("M122").Value = -19937.125
$ws.Range("N122").Value = -17865.727
$ws.Range("H132").Value = 6334.25
$ws.Range("I132").Value = 10124.5
$ws.Range("J132").Value = 3386.2778
$ws.Range("K132").Value = 30373.5
$ws.Range("L132").Value = 10158.8334
$ws.Range("M132").Value = -27843.5
$ws.Range("N132").Value = -15218.8334
$ws.Range("H136").Value = 2561.125
$ws.Range("I136").Value = 2095.8572
$ws.Range("J136").Value = 3449.3635
$ws.Range("K136").Value = 6287.571599999999
$ws.Range("L136").Value = 10348.0905
$ws.Range("M136").Value = -3737.571599999999
$ws.Range("N136").Value = -15448.0905
$ws.Range("H138").Value = 29800
$ws.Range("J138").Value = 29800
$ws.Range("L138").Value = 29800
$ws.Range("N138").Value = -40080
$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H3").Value = 24877.143
$ws.Range("I3").Value = 867.36365
$ws.Range("J3").Value = 51287.9
$ws.Range("K3").Value = 867.36365
$ws.Range("L3").Value = 51287.9
$ws.Range("M3").Value = -753.36365
$ws.Range("N3").Value = -51515.9
$ws.Range("H86").Value = 135627.33
$ws.Range("I86").Value = 3170.8572
$ws.Range("J86").Value = 251526.75
$ws.Range("K86").Value = 3170.8572
$ws.Range("L86").Value = 251526.75
$ws.Range("M86").Value = -2047.8572
$ws.Range("N86").Value = -253772.75
$ws.Range("H89").Value = 135627.33
$ws.Range("I89").Value = 3170.8572
$ws.Range("J89").Value = 251526.75
$ws.Range("K89").Value = 15854.286
$ws.Range("L89").Value = 1257633.75
$ws.Range("M89").Value = -10238.286
$ws.Range("N89").Value = -1268865.75
$ws.Range("H105").Value = 2335.647
$ws.Range("I105").Value = 2136.8572
$ws.Range("J105").Value = 3263.3333
$ws.Range("K105").Value = 2136.8572
$ws.Range("L105").Value = 3263.3333
$ws.Range("M105").Value = -389.8571999999999
$ws.Range("N105").Value = -6757.3333
$ws.Range("H134").Value = 1854.2
$ws.Range("I134").Value = 1599.262
$ws.Range("J134").Value = 2677.8462
$ws.Range("K134").Value = 4797.786
$ws.Range("L134").Value = 8033.5386
$ws.Range("M134").Value = -2262.786
$ws.Range("N134").Value = -13103.5386
$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H31").Value = 2308.4333
$ws.Range("I31").Value = 1358.5555
$ws.Range("J31").Value = 3733.25
$ws.Range("K31").Value = 1358.5555
$ws.Range("L31").Value = 3733.25
$ws.Range("M31").Value = -1063.5555
$ws.Range("N31").Value = -4323.25
$ws.Range("H34").Value = 2308.4333
$ws.Range("I34").Value = 1358.5555
$ws.Range("J34").Value = 3733.25
$ws.Range("K34").Value = 1358.5555
$ws.Range("L34").Value = 3733.25
$ws.Range("M34").Value = -1156.5555
$ws.Range("N34").Value = -4137.25
$ws.Range("H58").Value = 1853849.8
$ws.Range("I58").Value = 2471250.2
$ws.Range("J58").Value = 1647.8
$ws.Range("K58").Value = 2471250.2
$ws.Range("L58").Value = 1647.8
$ws.Range("M58").Value = -2471047.2
$ws.Range("N58").Value = -2053.8
$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 25000
$ws.Range("N125").Value = -29920
$ws.Range("H132").Value = 323512.2
$ws.Range("I132").Value = 467106.88
$ws.Range("J132").Value = 3185.6155
$ws.Range("K132").Value = 1401320.64
$ws.Range("L132").Value = 9556.8465
$ws.Range("M132").Value = -1398790.64
$ws.Range("N132").Value = -14616.8465
$ws.Range("H134").Value = 1828.9062
$ws.Range("I134").Value = 1387.3889
$ws.Range("J134").Value = 2396.5715
$ws.Range("K134").Value = 4162.1667
$ws.Range("L134").Value = 7189.7145
$ws.Range("M134").Value = -1627.1667
$ws.Range("N134").Value = -12259.7145
$ws.Range("H136").Value = 1853849.8
$ws.Range("I136").Value = 2471250.2
$ws.Range("J136").Value = 1647.8
$ws.Range("K136").Value = 7413750.600000001
$ws.Range("L136").Value = 4943.4
$ws.Range("M136").Value = -7411200.600000001
$ws.Range("N136").Value = -10043.4
$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H102").Value = 8626.1
$ws.Range("J102").Value = 8626.1
$ws.Range("L102").Value = 25878.3
$ws.Range("N102").Value = -30746.3
$ws.Range("H113").Value = 618.36365
$ws.Range("I113").Value = 618.36365
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1855.09095
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 314.90905
$ws.Range("N113").ClearContents()
$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H51").Value = 14914.833
$ws.Range("J51").Value = 14914.833
$ws.Range("L51").Value = 14914.833
$ws.Range("N51").Value = -15932.833
$ws.Range("H70").Value = 6973.4585
$ws.Range("I70").Value = 6568.143
$ws.Range("J70").Value = 7540.9
$ws.Range("K70").Value = 6568.143
$ws.Range("L70").Value = 7540.9
$ws.Range("M70").Value = -6298.143
$ws.Range("N70").Value = -8080.9
$ws.Range("H73").Value = 6973.4585
$ws.Range("I73").Value = 6568.143
$ws.Range("J73").Value = 7540.9
$ws.Range("K73").Value = 6568.143
$ws.Range("L73").Value = 7540.9
$ws.Range("M73").Value = -5632.143
$ws.Range("N73").Value = -9412.9
$ws.Range("H132").Value = 2297
$ws.Range("I132").Value = 1975.0209
$ws.Range("J132").Value = 3584.9167
$ws.Range("K132").Value = 5925.0627
$ws.Range("L132").Value = 10754.7501
$ws.Range("M132").Value = -3395.0627
$ws.Range("N132").Value = -15814.7501
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H40").Value = 3786.8823
$ws.Range("I40").Value = 3755.5
$ws.Range("J40").Value = 3933.3333
$ws.Range("K40").Value = 3755.5
$ws.Range("L40").Value = 3933.3333
$ws.Range("M40").Value = -3619.5
$ws.Range("N40").Value = -4205.3333
$ws.Range("H132").Value = 3130.366
$ws.Range("I132").Value = 2182.5833
$ws.Range("J132").Value = 4468.4116
$ws.Range("K132").Value = 6547.749899999999
$ws.Range("L132").Value = 13405.2348
$ws.Range("M132").Value = -4017.749899999999
$ws.Range("N132").Value = -18465.2348
$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H132").Value = 1861.1951
$ws.Range("I132").Value = 1197.25
$ws.Range("J132").Value = 4221.8887
$ws.Range("K132").Value = 3591.75
$ws.Range("L132").Value = 12665.6661
$ws.Range("M132").Value = -1061.75
$ws.Range("N132").Value = -17725.6661
$ws.Range("H136").Value = 1763.6923
$ws.Range("I136").Value = 1607.3182
$ws.Range("J136").Value = 2623.75
$ws.Range("K136").Value = 4821.9546
$ws.Range("L136").Value = 7871.25
$ws.Range("M136").Value = -2271.9546
$ws.Range("N136").Value = -12971.25
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
